$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.310.46"
$ws.Range("E2").Value = "  -1.57%  "

$ws.Range("D3").Value = "'1.551.36"
$ws.Range("E3").Value = "  -1.54%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'209.98"
$ws.Range("E5").Value = "  -1.70%  "

$ws.Range("E6").Value = "  -1.80%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").Value = "'23.75"
$ws.Range("E8").Value = "  -1.88%  "

$ws.Range("E9").Value = "  -1.86%  "

$ws.Range("E10").Value = "  -1.73%  "

$ws.Range("D11").Value = "'0.0891"
$ws.Range("E11").Value = "  +0.04%  "

$ws.Range("D12").Value = "'1.774.01"
$ws.Range("E12").Value = "  -1.52%  "

$ws.Range("D13").Value = "'1.551.42"
$ws.Range("E13").Value = "  -1.65%  "

$ws.Range("D14").Value = "'28.289.04"
$ws.Range("E14").Value = "  -1.65%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'3.62"
$ws.Range("E15").Value = "  -1.86%  "

$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.509"
$ws.Range("E16").Value = "  -2.63%  "

$ws.Range("D17").Value = "'60.56"
$ws.Range("E17").Value = "  -3.15%  "

$ws.Range("D18").Value = "'227.81"
$ws.Range("E18").Value = "  -1.76%  "

$ws.Range("E19").Value = "  -1.11%  "

$ws.Range("D20").Value = "'0.0₃0674"
$ws.Range("E20").Value = "  -2.82%  "

$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("E22").Value = "  +0.54%  "

$ws.Range("D23").Value = "'8.91"
$ws.Range("E23").Value = "  -3.00%  "

$ws.Range("E24").Value = "  -4.18%  "

$ws.Range("D25").Value = "'151.05"
$ws.Range("E25").Value = "  -0.97%  "

$ws.Range("E26").Value = "  -1.89%  "

$ws.Range("E27").Value = "  -1.65%  "

$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("D29").Value = "'6.24"
$ws.Range("E29").Value = "  -3.52%  "

$ws.Range("E30").Value = "  -3.48%  "

$ws.Range("E31").Value = "  -4.34%  "

$ws.Range("E32").Value = "  -1.65%  "

$ws.Range("D33").Value = "'1.386.92"
$ws.Range("E33").Value = "  -0.78%  "

$ws.Range("D34").Value = "'3.02"
$ws.Range("E34").Value = "  -3.37%  "

$ws.Range("E35").Value = "  +2.05%  "

$ws.Range("E36").Value = "  -4.45%  "

$ws.Range("E37").Value = "  -1.30%  "

$ws.Range("E38").Value = "  -1.35%  "

$ws.Range("E39").Value = "  -3.04%  "

$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "'0.512"
$ws.Range("E40").Value = "  -2.80%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'1.91"
$ws.Range("E41").Value = "  +0.61%  "

$ws.Range("E42").Value = "  -0.13%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'0.777"
$ws.Range("E43").Value = "  -2.28%  "

$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "'0.0464"
$ws.Range("E44").Value = "  -1.57%  "

$ws.Range("D45").Value = "'5.35"
$ws.Range("E45").Value = "  -3.01%  "

$ws.Range("D46").Value = "'61.77"
$ws.Range("E46").Value = "  -2.53%  "

$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'0.906"
$ws.Range("E47").Value = "  -5.99%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'85.60"
$ws.Range("E48").Value = "  -1.31%  "

$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").Value = "'42.85"
$ws.Range("E49").Value = "  +7.01%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.0₆0103"
$ws.Range("E50").Value = "  -0.07%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0508"
$ws.Range("E51").Value = "  -1.90%  "

